$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.180.10"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "1.627.31"
$ws.Range("E3").Value = "  -0.99%  "
$ws.Range("E4").Value = "  -0.37%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.78"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.522"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.55%  "
$ws.Range("E7").Value = "  -0.46%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.253"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0628"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.38"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.13%  "
$ws.Range("E11").Value = "  -0.15%  "
$ws.Range("D12").Value = "1.635.82"
$ws.Range("E12").Value = "  -0.67%  "
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("E14").Value = "  +0.24%  "
$ws.Range("D15").Value = "27.155.82"
$ws.Range("E15").Value = "  +0.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.70"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.82%  "
$ws.Range("E17").Value = "  +0.96%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "217.28"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.998"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.95"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.77%  "
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("E22").Value = "  -4.64%  "
$ws.Range("E23").Value = "  -1.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "147.58"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.998"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.53%  "
$ws.Range("E26").Value = "  -2.62%  "
$ws.Range("E27").Value = "  +0.12%  "
$ws.Range("E28").Value = "  -0.83%  "
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("E30").Value = "  -1.08%  "
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("E32").Value = "  -0.82%  "
$ws.Range("D33").Value = "1.351.26"
$ws.Range("E33").Value = "  +7.22%  "
$ws.Range("E34").Value = "  +0.23%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.45"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.76%  "
$ws.Range("E36").Value = "  -0.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.551"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.855"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.804"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.47%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "65.63"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +6.34%  "
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.23"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("E43").Value = "  -0.51%  "
$ws.Range("D44").Value = "1.765.74"
$ws.Range("E44").Value = "  -0.99%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "90.63"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.21%  "
$ws.Range("E46").Value = "  +1.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.846"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +27.51%  "
$ws.Range("E48").Value = "  +2.81%  "
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("D50").Value = "0.0₇0988"
$ws.Range("E50").Value = "  -7.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.60"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.35%  "
